$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Copy formatting (number format, style) from the last existing row (80) down to
# the new row (81) so the new row matches the look of the rest of the table.
$ws.Range("A80:E80").Copy()
$ws.Range("A81:E81").PasteSpecial(-4104)

# Fill in the new data for 31 May 2020 (serial date 43982)
$ws.Cells.Item(81, 1).Value = 43982
$ws.Cells.Item(81, 2).Value = 35600
$ws.Cells.Item(81, 3).Value = 839
$ws.Cells.Item(81, 4).Value = 25
$ws.Cells.Item(81, 5).Value = 989

# Grow the worksheet Table (ListObject) so it includes the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E81"))
